$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of an already-styled data row (row 5, style index 0)
# and apply it to row 14 (which previously lacked an explicit style) and to the
# new rows 15-28 that will hold the additional evaluation results.
$ws.Range("A5:F5").Copy()
$ws.Range("A14:F28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A15").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Right\noise_only.csv"
$ws.Range("B15").Value2 = 3.555207598209381
$ws.Range("C15").Value2 = 0.9829049305704268
$ws.Range("D15").Value2 = 33.32115936279297
$ws.Range("E15").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Right\test\clean"
$ws.Range("F15").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Right/noise_only"

$ws.Range("A16").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\noise_only.csv"
$ws.Range("B16").Value2 = 3.662087273597717
$ws.Range("C16").Value2 = 0.9853239165057607
$ws.Range("D16").Value2 = 33.90171051025391
$ws.Range("E16").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\test\clean"
$ws.Range("F16").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontRight/noise_only"

$ws.Range("A17").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\noise_only.csv"
$ws.Range("B17").Value2 = 2.170329254865647
$ws.Range("C17").Value2 = 0.903891390591925
$ws.Range("D17").Value2 = 14.7094898223877
$ws.Range("E17").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\test\clean"
$ws.Range("F17").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Front/noise_only"

$ws.Range("A18").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\noise_only.csv"
$ws.Range("B18").Value2 = 3.679636192321778
$ws.Range("C18").Value2 = 0.9864447880212298
$ws.Range("D18").Value2 = 33.60507965087891
$ws.Range("E18").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\test\clean"
$ws.Range("F18").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontLeft/noise_only"

$ws.Range("A19").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\noise_only.csv"
$ws.Range("B19").Value2 = 3.398727691173554
$ws.Range("C19").Value2 = 0.9825341801258716
$ws.Range("D19").Value2 = 31.5970344543457
$ws.Range("E19").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\test\clean"
$ws.Range("F19").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Left/noise_only"

$ws.Range("A20").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Right\noise_reverbe.csv"
$ws.Range("B20").Value2 = 2.303643715381622
$ws.Range("C20").Value2 = 0.942897948267561
$ws.Range("D20").Value2 = 19.03293228149414
$ws.Range("E20").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Right\test\clean"
$ws.Range("F20").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Right/noise_reverbe"

$ws.Range("A21").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\noise_reverbe.csv"
$ws.Range("B21").Value2 = 2.62207989692688
$ws.Range("C21").Value2 = 0.9608179637329928
$ws.Range("D21").Value2 = 21.1391773223877
$ws.Range("E21").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\test\clean"
$ws.Range("F21").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontRight/noise_reverbe"

$ws.Range("A22").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\reverbe_only.csv"
$ws.Range("B22").Value2 = 3.079035294055939
$ws.Range("C22").Value2 = 0.9805201789132024
$ws.Range("D22").Value2 = 21.89427185058594
$ws.Range("E22").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontRight\test\clean"
$ws.Range("F22").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontRight/reverbe_only"

$ws.Range("A23").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\noise_reverbe.csv"
$ws.Range("B23").Value2 = 1.692980843782425
$ws.Range("C23").Value2 = 0.907805576066866
$ws.Range("D23").Value2 = 14.99407958984375
$ws.Range("E23").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\test\clean"
$ws.Range("F23").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Front/noise_reverbe"

$ws.Range("A24").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\reverbe_only.csv"
$ws.Range("B24").Value2 = 2.923302006721497
$ws.Range("C24").Value2 = 0.9756118876948369
$ws.Range("D24").Value2 = 21.40922546386719
$ws.Range("E24").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Front\test\clean"
$ws.Range("F24").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Front/reverbe_only"

$ws.Range("A25").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\noise_reverbe.csv"
$ws.Range("B25").Value2 = 2.540942287445068
$ws.Range("C25").Value2 = 0.9593165930115036
$ws.Range("D25").Value2 = 20.86977577209473
$ws.Range("E25").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\test\clean"
$ws.Range("F25").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontLeft/noise_reverbe"

$ws.Range("A26").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\reverbe_only.csv"
$ws.Range("B26").Value2 = 3.053740119934082
$ws.Range("C26").Value2 = 0.9807686714588723
$ws.Range("D26").Value2 = 21.64970588684082
$ws.Range("E26").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\FrontLeft\test\clean"
$ws.Range("F26").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/FrontLeft/reverbe_only"

$ws.Range("A27").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\noise_reverbe.csv"
$ws.Range("B27").Value2 = 2.104562246799469
$ws.Range("C27").Value2 = 0.945898200075469
$ws.Range("D27").Value2 = 19.58468627929688
$ws.Range("E27").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\test\clean"
$ws.Range("F27").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Left/noise_reverbe"

$ws.Range("A28").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\evaluation\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\reverbe_only.csv"
$ws.Range("B28").Value2 = 2.900625514984131
$ws.Range("C28").Value2 = 0.9768721380945248
$ws.Range("D28").Value2 = 21.19450759887695
$ws.Range("E28").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\mix_data\subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm\Left\test\clean"
$ws.Range("F28").Value2 = "C:\Users\kataoka-lab\Desktop\sound_data\RESULT\output_wav/subset_DEMAND_hoth_1010dB_05sec_4ch_circular_10cm_Dtype/Left/reverbe_only"
